# Updates cryptos list prices (column D) and volume(1h) percentages (column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose new column-D price string looks like an ordinary decimal number
# (e.g. "568.88"). Without help these would be auto-converted from text to a
# number by Excel's input parser, unlike the original "dotted" price strings
# (e.g. "63.543.36") which are never valid numbers and stay text naturally.
$forceTextRows = @(5, 6, 9, 15, 19, 20, 21, 22, 24, 25, 27, 31, 32, 33, 34, 36, 37, 40, 41, 44, 45, 46, 47, 51)

$updates = @(
    @{ Row = 2;  D = "63.543.36";  E = "  +1.35%  " },
    @{ Row = 3;  D = "3.416.29";   E = "  +2.19%  " },
    @{ Row = 4;  D = $null;        E = "  +0.00%  " },
    @{ Row = 5;  D = "568.88";     E = "  +1.36%  " },
    @{ Row = 6;  D = "156.11";     E = "  +2.39%  " },
    @{ Row = 7;  D = $null;        E = "  -0.04%  " },
    @{ Row = 8;  D = "3.413.46";   E = "  +1.96%  " },
    @{ Row = 9;  D = "0.547";      E = "  +2.71%  " },
    @{ Row = 10; D = $null;        E = "  -0.12%  " },
    @{ Row = 11; D = $null;        E = "  +3.89%  " },
    @{ Row = 12; D = $null;        E = "  -0.37%  " },
    @{ Row = 13; D = "4.001.25";   E = "  +2.08%  " },
    @{ Row = 14; D = $null;        E = "  -2.92%  " },
    @{ Row = 15; D = "0.0000195";  E = "  +8.32%  " },
    @{ Row = 16; D = $null;        E = "  +1.16%  " },
    @{ Row = 17; D = "63.623.62";  E = "  +1.45%  " },
    @{ Row = 18; D = "3.407.97";   E = "  +2.52%  " },
    @{ Row = 19; D = "6.27";       E = "  -1.25%  " },
    @{ Row = 20; D = "14.12";      E = "  +2.28%  " },
    @{ Row = 21; D = "380.81";     E = "  -0.91%  " },
    @{ Row = 22; D = "8.09";       E = "  -3.39%  " },
    @{ Row = 23; D = $null;        E = "  +0.39%  " },
    @{ Row = 24; D = "71.74";      E = "  +2.31%  " },
    @{ Row = 25; D = "0.530";      E = "  -1.23%  " },
    @{ Row = 26; D = $null;        E = "  +27.95%  " },
    @{ Row = 27; D = "9.42";       E = "  +5.77%  " },
    @{ Row = 28; D = $null;        E = "  +0.13%  " },
    @{ Row = 29; D = $null;        E = "  +0.11%  " },
    @{ Row = 30; D = $null;        E = "  +8.47%  " },
    @{ Row = 31; D = "1.38";       E = "  +5.04%  " },
    @{ Row = 32; D = "2.01";       E = "  +1.27%  " },
    @{ Row = 33; D = "23.28";      E = "  +1.67%  " },
    @{ Row = 34; D = "6.40";       E = "  -2.26%  " },
    @{ Row = 35; D = $null;        E = "  +0.01%  " },
    @{ Row = 36; D = "6.81";       E = "  +1.64%  " },
    @{ Row = 37; D = "159.84";     E = "  -0.12%  " },
    @{ Row = 38; D = $null;        E = "  -1.21%  " },
    @{ Row = 39; D = "2.966.37";   E = "  +5.97%  " },
    @{ Row = 40; D = "0.0763";     E = "  +3.34%  " },
    @{ Row = 41; D = "27.10";      E = "  +0.61%  " },
    @{ Row = 42; D = $null;        E = "  -3.63%  " },
    @{ Row = 43; D = $null;        E = "  +0.96%  " },
    @{ Row = 44; D = "41.87";      E = "  +3.47%  " },
    @{ Row = 45; D = "0.764";      E = "  +2.57%  " },
    @{ Row = 46; D = "4.34";       E = "  +2.03%  " },
    @{ Row = 47; D = "23.33";      E = "  +6.09%  " },
    @{ Row = 48; D = $null;        E = "  +3.57%  " },
    @{ Row = 49; D = $null;        E = "  +23.65%  " },
    @{ Row = 50; D = $null;        E = "  +4.47%  " },
    @{ Row = 51; D = "6.37";       E = "  +0.92%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($forceTextRows -contains $r) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
